$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers for Wins / Losses / Ties in columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the existing header formatting (bold, centered, thin border) used by the rest of row 1
# by copying the format from an existing header cell rather than re-deriving it,
# so the new cells reuse the same style definition instead of creating a new one.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the team record (Wins/Losses/Ties) for every data row
for ($r = 2; $r -le 48; $r++) {
    $ws.Cells.Item($r, 30).Value = 77
    $ws.Cells.Item($r, 31).Value = 85
    $ws.Cells.Item($r, 32).Value = 0
}
